# Multiply the "size" columns (A, C, E, G, I, K) for data rows 2-19 by 1000.
# The "Count" columns (B, D, F, H, J, L) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "C", "E", "G", "I", "K")

for ($row = 2; $row -le 19; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($val -ne $null) {
            $cell.Value2 = $val * 1000
        }
    }
}
